$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.817.71"
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").Value = "3.087.24"
$ws.Range("E3").Value = "  +4.96%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.75"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.43"
$ws.Range("E6").Value = "  +5.97%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.081.91"
$ws.Range("E8").Value = "  +4.92%  "
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.62"
$ws.Range("E10").Value = "  -1.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +3.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.483"
$ws.Range("E12").Value = "  +5.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.44"
$ws.Range("E14").Value = "  +6.17%  "
$ws.Range("D16").Value = "3.599.97"
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("D17").Value = "66.826.24"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.18"
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "3.087.53"
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.13"
$ws.Range("E20").Value = "  +8.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.27"
$ws.Range("E21").Value = "  +4.65%  "
$ws.Range("E22").Value = "  +4.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("E23").Value = "  +3.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.41"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("E25").Value = "  +6.69%  "
$ws.Range("E26").Value = "  +7.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("E31").Value = "  +3.70%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.23"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("E34").Value = "  +3.53%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.88"
$ws.Range("E37").Value = "  +3.09%  "
$ws.Range("E38").Value = "  +6.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "46.98"
$ws.Range("E39").Value = "  +5.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.23"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.317"
$ws.Range("E41").Value = "  +6.58%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  +2.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "382.61"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").Value = "2.772.58"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.19"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.79"
$ws.Range("E50").Value = "  +6.59%  "
$ws.Range("E51").Value = "  +1.56%  "
